# Auto-generated edit script: updates Leve profit-tracking values
# across 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the
# scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 500.5
$ws.Range("I33").Value = 572.25
$ws.Range("J33").Value = 357
$ws.Range("K33").Value = 572.25
$ws.Range("L33").Value = 357
$ws.Range("M33").Value = -343.25
$ws.Range("N33").Value = -815

$ws.Range("H40").Value = 4847.6313
$ws.Range("I40").Value = 4374.5
$ws.Range("J40").Value = 5191.727
$ws.Range("K40").Value = 4374.5
$ws.Range("L40").Value = 5191.727
$ws.Range("M40").Value = -4199.5
$ws.Range("N40").Value = -5541.727

$ws.Range("H106").Value = 704.6
$ws.Range("I106").Value = 704.6
$ws.Range("K106").Value = 704.6
$ws.Range("M106").Value = -73.60000000000002

$ws.Range("H132").Value = 8751.290000000001
$ws.Range("I132").Value = 7975.552
$ws.Range("K132").Value = 23926.656
$ws.Range("M132").Value = -21396.656

$ws.Range("H135").Value = 570.5
$ws.Range("I135").Value = 175.85715
$ws.Range("J135").Value = 3333
$ws.Range("K135").Value = 1582.71435
$ws.Range("L135").Value = 29997
$ws.Range("M135").Value = 952.28565
$ws.Range("N135").Value = -35067

$ws.Range("H138").Value = 2342.9092
$ws.Range("I138").Value = 1917.5555
$ws.Range("J138").Value = 4257
$ws.Range("K138").Value = 5752.666499999999
$ws.Range("L138").Value = 12771
$ws.Range("M138").Value = -612.6664999999994
$ws.Range("N138").Value = -23051

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 16000

$ws.Range("H74").Value = 1602.3636
$ws.Range("I74").Value = 1408
$ws.Range("J74").Value = 2833.3333
$ws.Range("K74").Value = 1408
$ws.Range("L74").Value = 2833.3333
$ws.Range("M74").Value = -534
$ws.Range("N74").Value = -4581.3333

$ws.Range("H77").Value = 1602.3636
$ws.Range("I77").Value = 1408
$ws.Range("J77").Value = 2833.3333
$ws.Range("K77").Value = 7040
$ws.Range("L77").Value = 14166.6665
$ws.Range("M77").Value = -2672
$ws.Range("N77").Value = -22902.6665

$ws.Range("H102").Value = 7181.385
$ws.Range("I102").Value = 3073.8
$ws.Range("J102").Value = 9748.625
$ws.Range("K102").Value = 3073.8
$ws.Range("L102").Value = 9748.625
$ws.Range("M102").Value = -1451.8
$ws.Range("N102").Value = -12992.625

$ws.Range("H110").Value = 2684.3872
$ws.Range("I110").Value = 1234.2941
$ws.Range("K110").Value = 1234.2941
$ws.Range("M110").Value = 810.7058999999999

$ws.Range("H122").Value = 1262.5
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H132").Value = 6193.1665
$ws.Range("I132").Value = 6431.8
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 19295.4
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -16765.4
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 475
$ws.Range("I12").Value = 475
$ws.Range("K12").Value = 475
$ws.Range("M12").Value = -307

$ws.Range("H20").Value = 4171.875
$ws.Range("I20").Value = 4483.5713
$ws.Range("K20").Value = 4483.5713
$ws.Range("M20").Value = -4236.5713

$ws.Range("H63").Value = 84999.5
$ws.Range("J63").Value = 84999.5
$ws.Range("L63").Value = 84999.5
$ws.Range("N63").Value = -86371.5

$ws.Range("H66").Value = 84999.5
$ws.Range("J66").Value = 84999.5
$ws.Range("L66").Value = 254998.5
$ws.Range("N66").Value = -261862.5

$ws.Range("H107").Value = 3078.5293
$ws.Range("I107").Value = 1357.3478
$ws.Range("J107").Value = 6677.364
$ws.Range("K107").Value = 1357.3478
$ws.Range("L107").Value = 6677.364
$ws.Range("M107").Value = 562.6522
$ws.Range("N107").Value = -10517.364

$ws.Range("H134").Value = 5110
$ws.Range("I134").Value = 5110
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 15330
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -12795
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2359.7
$ws.Range("I16").Value = 1999.5714
$ws.Range("K16").Value = 1999.5714
$ws.Range("M16").Value = -1712.5714

$ws.Range("H99").Value = 5018.3335
$ws.Range("I99").Value = 4547.1
$ws.Range("K99").Value = 4547.1
$ws.Range("M99").Value = -3049.1

$ws.Range("H107").Value = 619.8
$ws.Range("J107").Value = 714.25
$ws.Range("L107").Value = 714.25
$ws.Range("N107").Value = -4554.25

$ws.Range("H113").Value = 2359.7
$ws.Range("I113").Value = 1999.5714
$ws.Range("K113").Value = 1999.5714
$ws.Range("M113").Value = 170.4286

$ws.Range("H117").Value = 39979
$ws.Range("J117").Value = 39979
$ws.Range("L117").Value = 39979
$ws.Range("N117").Value = -49157

$ws.Range("H126").Value = 5018.3335
$ws.Range("I126").Value = 4547.1
$ws.Range("K126").Value = 13641.3
$ws.Range("M126").Value = -11171.3

$ws.Range("H132").Value = 1852.8636
$ws.Range("I132").Value = 1528.4
$ws.Range("J132").Value = 5097.5
$ws.Range("K132").Value = 4585.200000000001
$ws.Range("L132").Value = 15292.5
$ws.Range("M132").Value = -2055.200000000001
$ws.Range("N132").Value = -20352.5

$ws.Range("H134").Value = 2374.2354
$ws.Range("I134").Value = 1917.6451
$ws.Range("K134").Value = 5752.9353
$ws.Range("M134").Value = -3217.9353

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5700.25
$ws.Range("I80").Value = 1397.5
$ws.Range("K80").Value = 1397.5
$ws.Range("M80").Value = -399.5

$ws.Range("H83").Value = 5700.25
$ws.Range("I83").Value = 1397.5
$ws.Range("K83").Value = 6987.5
$ws.Range("M83").Value = -1995.5

$ws.Range("H113").Value = 7416.0557
$ws.Range("I113").Value = 5958.9
$ws.Range("K113").Value = 5958.9
$ws.Range("M113").Value = -3788.9

$ws.Range("H132").Value = 265755.5
$ws.Range("I132").Value = 349341
$ws.Range("J132").Value = 14999
$ws.Range("K132").Value = 1048023
$ws.Range("L132").Value = 44997
$ws.Range("M132").Value = -1045493
$ws.Range("N132").Value = -50057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2725.7144
$ws.Range("I2").Value = 2000
$ws.Range("J2").Value = 2846.6667
$ws.Range("K2").Value = 2000
$ws.Range("L2").Value = 2846.6667
$ws.Range("M2").Value = -1888
$ws.Range("N2").Value = -3070.6667

$ws.Range("H16").Value = 991.1667
$ws.Range("I16").Value = 723.5
$ws.Range("J16").Value = 1125
$ws.Range("K16").Value = 723.5
$ws.Range("L16").Value = 1125
$ws.Range("M16").Value = -553.5
$ws.Range("N16").Value = -1465

$ws.Range("H61").Value = 4319
$ws.Range("I61").Value = 3092.0667
$ws.Range("K61").Value = 3092.0667
$ws.Range("M61").Value = -2890.0667

$ws.Range("H76").Value = 20288
$ws.Range("J76").Value = 20288
$ws.Range("L76").Value = 20288
$ws.Range("N76").Value = -20964

$ws.Range("H79").Value = 20288
$ws.Range("J79").Value = 20288
$ws.Range("L79").Value = 20288
$ws.Range("N79").Value = -22628

$ws.Range("H113").Value = 4319
$ws.Range("I113").Value = 3092.0667
$ws.Range("K113").Value = 3092.0667
$ws.Range("M113").Value = -922.0666999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 7629.4
$ws.Range("I4").Value = 8428.416999999999
$ws.Range("J4").Value = 4433.3335
$ws.Range("K4").Value = 8428.416999999999
$ws.Range("L4").Value = 4433.3335
$ws.Range("M4").Value = -8315.416999999999
$ws.Range("N4").Value = -4659.3335

$ws.Range("H107").Value = 610.95
$ws.Range("I107").Value = 610
$ws.Range("J107").Value = 613.8
$ws.Range("K107").Value = 1830
$ws.Range("L107").Value = 1841.4
$ws.Range("M107").Value = 90
$ws.Range("N107").Value = -5681.4

$ws.Range("H116").Value = 1000000
$ws.Range("J116").Value = 1000000
$ws.Range("L116").Value = 1000000
$ws.Range("N116").Value = -1009178

$ws.Range("H122").Value = 2995.1516
$ws.Range("I122").Value = 2059.3
$ws.Range("K122").Value = 6177.900000000001
$ws.Range("M122").Value = -3727.900000000001

$ws.Range("H126").Value = 3614.1428
$ws.Range("I126").Value = 2357
$ws.Range("J126").Value = 6128.4287
$ws.Range("K126").Value = 7071
$ws.Range("L126").Value = 18385.2861
$ws.Range("M126").Value = -4601
$ws.Range("N126").Value = -23325.2861

$ws.Range("H132").Value = 3041.2046
$ws.Range("I132").Value = 2806.1428
$ws.Range("K132").Value = 8418.428400000001
$ws.Range("M132").Value = -5888.428400000001

$ws.Range("H136").Value = 3779.3262
$ws.Range("I136").Value = 2424.9429
$ws.Range("K136").Value = 7274.8287
$ws.Range("M136").Value = -4724.8287
